$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restyle column A (rows 2,4-12) to match the style already used by A3 ---
# (A3 uses cell-style index 12 - centred, default font/no explicit border override;
#  the other rows in column A still use the older style index 1.)
[void]$ws.Range("A3").Copy()
[void]$ws.Range("A2").PasteSpecial(-4122)
[void]$ws.Range("A4:A12").PasteSpecial(-4122)

# --- Remove the explicit left-aligned style override from several F cells ---
# (F2/F3 already carry no explicit style; copy that "no style" formatting onto
#  the other F cells that currently still have the old explicit style applied.)
[void]$ws.Range("F3").Copy()
[void]$ws.Range("F4").PasteSpecial(-4122)
[void]$ws.Range("F5").PasteSpecial(-4122)
[void]$ws.Range("F6").PasteSpecial(-4122)
[void]$ws.Range("F7").PasteSpecial(-4122)
[void]$ws.Range("F11").PasteSpecial(-4122)
[void]$ws.Range("F12").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Update the "ActualResult" text for a few rows to the shorter message ---
$ws.Range("F2").Value = "No message found"
$ws.Range("F9").Value = "No message found"
$ws.Range("F12").Value = "No message found"

# --- Move the active selection from F2 to A2 ---
[void]$ws.Range("A2").Select()
